# Add the new "AttenuationCorrection" column (AA) to the metadata header row,
# mirroring the style/formatting of the existing header cells and updating
# the column width + selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell with the same green header style used by the rest of row 1.
$ws.Range("AA1").Value = "AttenuationCorrection"
$ws.Range("AA1").Interior.Color = $ws.Range("Z1").Interior.Color

# Give the new column a sensible custom width (matches the other header cols,
# target stored width ~20.332).
$ws.Columns("AA").ColumnWidth = 19.5

# Move the selection to the newly added cell (mirrors the scrolled view in
# the saved workbook, where the new column is now in focus).
$ws.Range("T1").Select()
$ws.Range("AA1").Select()
